$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 26
$ws.Range("H26").Value = 400
$ws.Range("I26").Value = 400
$ws.Range("K26").Value = 400
$ws.Range("M26").Value = -56
# Row 100
$ws.Range("H100").Value = 1485.091
$ws.Range("I100").Value = 1333.6
$ws.Range("K100").Value = 1333.6
$ws.Range("M100").Value = -792.5999999999999
# Row 140
$ws.Range("H140").Value = 81914.5
$ws.Range("J140").Value = 81914.5
$ws.Range("L140").Value = 81914.5
$ws.Range("N140").Value = -92274.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 8141.367
$ws.Range("I61").Value = 4359.2144
$ws.Range("J61").Value = 13184.238
$ws.Range("K61").Value = 4359.2144
$ws.Range("L61").Value = 13184.238
$ws.Range("M61").Value = -4147.2144
$ws.Range("N61").Value = -13608.238
# Row 97
$ws.Range("H97").Value = 1308.8823
$ws.Range("I97").Value = 1203.6364
$ws.Range("K97").Value = 1203.6364
$ws.Range("M97").Value = -707.6364000000001
# Row 121
$ws.Range("H121:L121").ClearContents()
$ws.Range("N121").ClearContents()
# Row 122
$ws.Range("H122:N122").ClearContents()
# Row 123
$ws.Range("H123:L123").ClearContents()
$ws.Range("N123").ClearContents()
# Row 124
$ws.Range("H124:L124").ClearContents()
# Row 125
$ws.Range("H125:L125").ClearContents()
$ws.Range("N125").ClearContents()
# Row 126
$ws.Range("H126:M126").ClearContents()
# Row 127
$ws.Range("H127:L127").ClearContents()
$ws.Range("N127").ClearContents()
# Row 128
$ws.Range("H128:L128").ClearContents()
$ws.Range("N128").ClearContents()
# Row 129
$ws.Range("H129:L129").ClearContents()
$ws.Range("N129").ClearContents()
# Row 130
$ws.Range("H130:L130").ClearContents()
$ws.Range("N130").ClearContents()
# Row 131
$ws.Range("H131:L131").ClearContents()
# Row 132
$ws.Range("H132:N132").ClearContents()
# Row 133
$ws.Range("H133:L133").ClearContents()
$ws.Range("N133").ClearContents()
# Row 134
$ws.Range("H134:L134").ClearContents()
$ws.Range("N134").ClearContents()
# Row 135
$ws.Range("H135:L135").ClearContents()
$ws.Range("N135").ClearContents()
# Row 136
$ws.Range("H136").Value = 8141.367
$ws.Range("I136").Value = 4359.2144
$ws.Range("J136").Value = 13184.238
$ws.Range("K136").Value = 13077.6432
$ws.Range("L136").Value = 39552.714
$ws.Range("M136").Value = -10527.6432
$ws.Range("N136").Value = -44652.714
# Row 137
$ws.Range("H137:L137").ClearContents()
# Row 138
$ws.Range("H138:L138").ClearContents()
$ws.Range("N138").ClearContents()
# Row 139
$ws.Range("H139:L139").ClearContents()
$ws.Range("N139").ClearContents()
# Row 140
$ws.Range("H140:L140").ClearContents()
$ws.Range("N140").ClearContents()
# Row 141
$ws.Range("H141:L141").ClearContents()
$ws.Range("N141").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 1176.5555
$ws.Range("I94").Value = 1084.1428
$ws.Range("K94").Value = 1084.1428
$ws.Range("M94").Value = -633.1428000000001
# Row 99
$ws.Range("H99").Value = 1069.7693
$ws.Range("I99").Value = 1069.7693
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1069.7693
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 428.2307000000001
$ws.Range("N99").ClearContents()
# Row 134
$ws.Range("H134").Value = 31151.572
$ws.Range("I134").Value = 2519.0952
$ws.Range("J134").Value = 74100.28999999999
$ws.Range("K134").Value = 7557.285600000001
$ws.Range("L134").Value = 222300.87
$ws.Range("M134").Value = -5022.285600000001
$ws.Range("N134").Value = -227370.87

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 134
$ws.Range("H134").Value = 2466.9524
$ws.Range("I134").Value = 2091.838
$ws.Range("J134").Value = 5242.8
$ws.Range("K134").Value = 6275.514000000001
$ws.Range("L134").Value = 15728.4
$ws.Range("M134").Value = -3740.514000000001
$ws.Range("N134").Value = -20798.4

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 122
$ws.Range("H122").Value = 1121.15
$ws.Range("I122").Value = 838
$ws.Range("J122").Value = 1161.6
$ws.Range("K122").Value = 7542
$ws.Range("L122").Value = 10454.4
$ws.Range("M122").Value = -5092
$ws.Range("N122").Value = -15354.4

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 2001
$ws.Range("I97").Value = 2001.1428
$ws.Range("K97").Value = 2001.1428
$ws.Range("M97").Value = -1505.1428
# Row 132
$ws.Range("H132").Value = 7873.231
$ws.Range("I132").Value = 8640.0625
$ws.Range("J132").Value = 7339.7827
$ws.Range("K132").Value = 25920.1875
$ws.Range("L132").Value = 22019.3481
$ws.Range("M132").Value = -23390.1875
$ws.Range("N132").Value = -27079.3481

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 93
$ws.Range("H93").Value = 2073.4443
$ws.Range("I93").Value = 1302.75
$ws.Range("J93").Value = 2690
$ws.Range("K93").Value = 1302.75
$ws.Range("L93").Value = 2690
$ws.Range("M93").Value = -54.75
$ws.Range("N93").Value = -5186
# Row 132
$ws.Range("H132").Value = 3485
$ws.Range("I132").Value = 3365.7693
$ws.Range("J132").Value = 3829.4443
$ws.Range("K132").Value = 10097.3079
$ws.Range("L132").Value = 11488.3329
$ws.Range("M132").Value = -7567.3079
$ws.Range("N132").Value = -16548.3329

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 5749
$ws.Range("I96").Value = 1498
$ws.Range("J96").Value = 10000
$ws.Range("K96").Value = 1498
$ws.Range("L96").Value = 10000
$ws.Range("M96").Value = -125
$ws.Range("N96").Value = -12746
# Row 123
$ws.Range("H123").Value = 53664.5
$ws.Range("J123").Value = 53664.5
$ws.Range("L123").Value = 53664.5
$ws.Range("N123").Value = -63464.5
# Row 132
$ws.Range("H132").Value = 2360.158
$ws.Range("I132").Value = 2530.3794
$ws.Range("J132").Value = 1811.6666
$ws.Range("K132").Value = 7591.138199999999
$ws.Range("L132").Value = 5434.9998
$ws.Range("M132").Value = -5061.138199999999
$ws.Range("N132").Value = -10494.9998
